$d = $word.ActiveDocument

# --- 1) "Groups are defined within the sample name column." -> "Let column name field blank." ---
$d.Content.Find.Execute(
    "Groups are defined within the sample name column.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Let column name field blank.", 2) | Out-Null

# Nudge formatting on the new text so Word keeps an explicit (empty) run-properties
# element on the run, matching the rest of the document's run style.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Let column name field blank.`r") {
        $r = $p.Range
        $r.MoveEnd(1, -1) | Out-Null
        $r.Font.Bold = 1
        $r.Font.Bold = 0
        break
    }
}

# --- 2) "order :" -> "Groups" + " :" (as two separate runs), scoped to its own paragraph ---
# Find the paragraph that now reads "order :IPSC,NPC,DA4W,DA6W"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("order :IPSC,NPC,DA4W,DA6W")) {
        $target = $p
        break
    }
}

$rng = $target.Range
$rng.Find.Execute("order :", $true, $false, $false, $false, $false, $true, 1, $false, "Groups :", 2) | Out-Null

# Re-fetch the paragraph range (text changed) and split "Groups :" into two runs by
# nudging formatting on each half - this forces Word to keep them as distinct <w:r> runs.
$full = $target.Range
$groupsPart = $d.Range($full.Start, $full.Start + 6)
$groupsPart.Font.Bold = 1
$groupsPart.Font.Bold = 0

$colonPart = $d.Range($full.Start + 6, $full.Start + 8)
$colonPart.Font.Bold = 1
$colonPart.Font.Bold = 0

# --- 3) Bookmarks: add a new __DdeLink__302_1425697230 bookmark (id 0) wrapping
#        "IPSC,NPC,DA4W,DA6W", and keep __DdeLink__12249_2669968041 as the outer one (id 1). ---
$valueRange = $target.Range
$valueRange.Find.Execute("IPSC,NPC,DA4W,DA6W", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmStart = $valueRange.Start
$bmEnd = $valueRange.End

$d.Bookmarks("__DdeLink__12249_2669968041").Delete()

$newBmRange = $d.Range($bmStart, $bmEnd)
$newBmRange.Bookmarks.Add("__DdeLink__302_1425697230") | Out-Null

$oldBmRange = $d.Range($bmStart, $bmEnd)
$oldBmRange.Bookmarks.Add("__DdeLink__12249_2669968041") | Out-Null
